$wb = $excel.ActiveWorkbook

# --- Rename the first sheet "En" -> "Entities" ---
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws1.Name = "Entities"

# --- Insert a new "name" column (B) into both sheets, shifting existing
#     columns one to the right, without touching the <cols> definitions. ---
foreach ($ws in @($ws1, $ws2)) {

    # Shift header row (row 1): columns K..B -> L..C
    for ($c = 11; $c -ge 2; $c--) {
        $srcCell = $ws.Cells.Item(1, $c)
        $dstCell = $ws.Cells.Item(1, $c + 1)
        $dstCell.Value = $srcCell.Value2
    }

    # Shift data row (row 2): columns E..B -> F..C, carrying styles along
    for ($c = 5; $c -ge 2; $c--) {
        $srcCell = $ws.Cells.Item(2, $c)
        $dstCell = $ws.Cells.Item(2, $c + 1)
        $dstCell.Value = $srcCell.Value2
        if ($c -ge 3) {
            $dstCell.Style = "Good"
            $dstCell.WrapText = $true
        } else {
            $dstCell.Style = "Good"
        }
    }

    # New column B content
    $ws.Cells.Item(1, 2).Value = "name"
    $ws.Cells.Item(2, 2).Value = "Public_Phone_Maintenance_Book"
    $ws.Cells.Item(2, 2).Style = "Good"

    # Row 2 grows to the max row height because the long, wrapped entries
    # now sit in much narrower columns than before.
    $ws.Rows.Item(2).RowHeight = 409.6
}

# --- Selection / active tab bookkeeping ---
# Sheet2 keeps a plain selection at B2 and is no longer the active tab.
$ws2.Range("B2").Select() | Out-Null
# Sheet1 ("Entities") becomes the active tab, keeping its original D8 selection.
$ws1.Range("D8").Select() | Out-Null
